# Adds three new year columns (2019, 2020, 2021) to the table, extending
# it from column P to column R, matching the existing table formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (divider row with bottom border) ---
# P3 currently carries the "plain" style; make it consistent with the
# rest of the divider row, and extend the border into the two new cells.
$ws.Range("P3:R3").Style = $ws.Range("O3").Style

# --- Row 4 (bold year headers, bordered) ---
$ws.Range("P4").Value = 2019
$ws.Range("Q4").Value = 2020
$ws.Range("R4").Value = 2021

$ws.Range("P4:R4").Font.Name = "Times New Roman"
$ws.Range("P4:R4").Font.Size = 10
$ws.Range("P4:R4").Font.Bold = $true
$ws.Range("P4:R4").Borders.Item(9).LineStyle = 6   # xlContinuous
$ws.Range("P4:R4").Borders.Item(9).Weight = -4138  # xlMedium
$ws.Range("P4:R4").VerticalAlignment = -4108       # xlCenter

# --- Row 5 (data values row, bordered top+bottom) ---
$ws.Range("P5").Value = 12.9
$ws.Range("Q5").Value = 15.2
$ws.Range("R5").Value = 10.4

$ws.Range("P5:R5").Font.Name = "Times New Roman"
$ws.Range("P5:R5").Font.Size = 10
$ws.Range("P5:R5").Font.Bold = $false
$ws.Range("P5:R5").Borders.Item(8).LineStyle = 6   # xlContinuous
$ws.Range("P5:R5").Borders.Item(8).Weight = -4138  # xlMedium
$ws.Range("P5:R5").Borders.Item(9).LineStyle = 6   # xlContinuous
$ws.Range("P5:R5").Borders.Item(9).Weight = -4138  # xlMedium
$ws.Range("P5:R5").VerticalAlignment = -4108       # xlCenter

# --- Selection / view state ---
$ws.Range("S3").Select()
